$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2 through 68 -> 7295
$ws.Range("C2:C68").Value = 7295

# Rows 69 through 252 -> 7293
$ws.Range("C69:C252").Value = 7293
